# update code tinh luong cho Quyen
# Bump last_edited_time for Thang 7 rows (7-12) to 2024-07-19T12:51:00.000Z
# and refresh the recalculated formula numbers on row 7 (Thang 7 page).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# last_edited_time (column D) for rows 7-12
$ws.Range("D7").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("D8").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("D9").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("D10").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("D11").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("D12").Value = "2024-07-19T12:51:00.000Z"

# Row 7 (Thang 7) recalculated formula numbers
$ws.Range("W7").Value = 217228000   # properties.Chi tieu.number
$ws.Range("AA7").Value = 171590000  # properties.Luy ke.formula.number
$ws.Range("AE7").Value = 388818000  # properties.Tong doanh thu.formula.number
$ws.Range("AN7").Value = 63000000   # properties.Thu no.number
